# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh to Sheets/Spriggan_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7386.0713
$ws.Range("I9").Value = 12721
$ws.Range("K9").Value = 12721
$ws.Range("M9").Value = -12552
$ws.Range("H15").Value = 1250487
$ws.Range("I15").Value = 1250487
$ws.Range("K15").Value = 3751461
$ws.Range("M15").Value = -3751292
$ws.Range("H17").Value = 441402.62
$ws.Range("J17").Value = 472905.1
$ws.Range("L17").Value = 1418715.3
$ws.Range("N17").Value = -1419051.3
$ws.Range("H43").Value = 6125
$ws.Range("I43").Value = 6125
$ws.Range("K43").Value = 6125
$ws.Range("M43").Value = -6056
$ws.Range("H106").Value = 2616.75
$ws.Range("I106").Value = 2722.3333
$ws.Range("J106").Value = 2300
$ws.Range("K106").Value = 2722.3333
$ws.Range("L106").Value = 2300
$ws.Range("M106").Value = -2091.3333
$ws.Range("N106").Value = -3562
$ws.Range("H132").Value = 3183.3076
$ws.Range("I132").Value = 3248.5833
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 9745.749899999999
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -7215.749899999999
$ws.Range("N132").Value = -12260
$ws.Range("H141").Value = 4099.75
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7071.091
$ws.Range("I32").Value = 5979.2383
$ws.Range("K32").Value = 5979.2383
$ws.Range("M32").Value = -5692.2383
$ws.Range("H97").Value = 883
$ws.Range("I97").Value = 999.5714
$ws.Range("J97").Value = 339
$ws.Range("K97").Value = 999.5714
$ws.Range("L97").Value = 339
$ws.Range("M97").Value = -503.5714
$ws.Range("N97").Value = -1331
$ws.Range("H132").Value = 10002283
$ws.Range("I132").Value = 10002283
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30006849
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -30004319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 79999
$ws.Range("J13").Value = 79999
$ws.Range("L13").Value = 79999
$ws.Range("N13").Value = -80335
$ws.Range("H20").Value = 1332.7273
$ws.Range("I20").Value = 993.5333000000001
$ws.Range("J20").Value = 2059.5715
$ws.Range("K20").Value = 993.5333000000001
$ws.Range("L20").Value = 2059.5715
$ws.Range("M20").Value = -746.5333000000001
$ws.Range("N20").Value = -2553.5715
$ws.Range("H80").Value = 663.7778
$ws.Range("I80").Value = 496.6
$ws.Range("J80").Value = 872.75
$ws.Range("K80").Value = 496.6
$ws.Range("L80").Value = 872.75
$ws.Range("M80").Value = 501.4
$ws.Range("N80").Value = -2868.75
$ws.Range("H83").Value = 663.7778
$ws.Range("I83").Value = 496.6
$ws.Range("J83").Value = 872.75
$ws.Range("K83").Value = 2483
$ws.Range("L83").Value = 4363.75
$ws.Range("M83").Value = 2509
$ws.Range("N83").Value = -14347.75
$ws.Range("H94").Value = 3350.7827
$ws.Range("I94").Value = 3537.8667
$ws.Range("K94").Value = 3537.8667
$ws.Range("M94").Value = -3086.8667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 52996.332
$ws.Range("J20").Value = 52996.332
$ws.Range("L20").Value = 52996.332
$ws.Range("N20").Value = -53468.332
$ws.Range("H30").Value = 52996.332
$ws.Range("J30").Value = 52996.332
$ws.Range("L30").Value = 52996.332
$ws.Range("N30").Value = -53178.332
$ws.Range("H105").Value = 974.0714
$ws.Range("I105").Value = 1035.8
$ws.Range("J105").Value = 819.75
$ws.Range("K105").Value = 1035.8
$ws.Range("L105").Value = 819.75
$ws.Range("M105").Value = 711.2
$ws.Range("N105").Value = -4313.75
$ws.Range("H121").Value = 89460.60000000001
$ws.Range("J121").Value = 89460.60000000001
$ws.Range("L121").Value = 89460.60000000001
$ws.Range("N121").Value = -92080.60000000001
$ws.Range("H122").Value = 2259.9
$ws.Range("I122").Value = 2259.9
$ws.Range("K122").Value = 6779.700000000001
$ws.Range("M122").Value = -4329.700000000001
$ws.Range("H128").Value = 52996.332
$ws.Range("J128").Value = 52996.332
$ws.Range("L128").Value = 52996.332
$ws.Range("N128").Value = -62956.332
$ws.Range("H134").Value = 50003140
$ws.Range("I134").Value = 62501424
$ws.Range("K134").Value = 187504272
$ws.Range("M134").Value = -187501737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 287.9091
$ws.Range("I14").Value = 287.9091
$ws.Range("K14").Value = 863.7273
$ws.Range("M14").Value = -690.7273
$ws.Range("H21").Value = 1030.6666
$ws.Range("J21").Value = 1501
$ws.Range("L21").Value = 4503
$ws.Range("N21").Value = -4849
$ws.Range("H34").Value = 663.3333
$ws.Range("I34").Value = 732
$ws.Range("K34").Value = 2196
$ws.Range("M34").Value = -2112
$ws.Range("H94").Value = 12313.909
$ws.Range("I94").Value = 2700
$ws.Range("J94").Value = 17807.572
$ws.Range("K94").Value = 8100
$ws.Range("L94").Value = 53422.716
$ws.Range("M94").Value = -7424
$ws.Range("N94").Value = -54774.716
$ws.Range("H139").Value = 4514.5
$ws.Range("I139").Value = 4514.5
$ws.Range("K139").Value = 13543.5
$ws.Range("M139").Value = -8403.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76372
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231864
$ws.Range("H107").Value = 5430.8335
$ws.Range("I107").Value = 4516.4
$ws.Range("J107").Value = 10003
$ws.Range("K107").Value = 4516.4
$ws.Range("L107").Value = 10003
$ws.Range("M107").Value = -2596.4
$ws.Range("N107").Value = -13843
$ws.Range("H122").Value = 2357.394
$ws.Range("I122").Value = 1271.375
$ws.Range("J122").Value = 5253.4443
$ws.Range("K122").Value = 3814.125
$ws.Range("L122").Value = 15760.3329
$ws.Range("M122").Value = -1364.125
$ws.Range("N122").Value = -20660.3329
$ws.Range("H132").Value = 6581176.5
$ws.Range("I132").Value = 8335357.5
$ws.Range("J132").Value = 2999.25
$ws.Range("K132").Value = 25006072.5
$ws.Range("L132").Value = 8997.75
$ws.Range("M132").Value = -25003542.5
$ws.Range("N132").Value = -14057.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2762.5
$ws.Range("I7").Value = 2863.08
$ws.Range("J7").Value = 2259.6
$ws.Range("K7").Value = 2863.08
$ws.Range("L7").Value = 2259.6
$ws.Range("M7").Value = -2751.08
$ws.Range("N7").Value = -2483.6
$ws.Range("H61").Value = 3144.1052
$ws.Range("I61").Value = 3144.1052
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3144.1052
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2942.1052
$ws.Range("H68").Value = 2874.75
$ws.Range("I68").Value = 2833
$ws.Range("K68").Value = 2833
$ws.Range("M68").Value = -2084
$ws.Range("H71").Value = 2874.75
$ws.Range("I71").Value = 2833
$ws.Range("K71").Value = 14165
$ws.Range("M71").Value = -10421
$ws.Range("H113").Value = 3144.1052
$ws.Range("I113").Value = 3144.1052
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3144.1052
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -974.1052
$ws.Range("H118").Value = 79999
$ws.Range("J118").Value = 79999
$ws.Range("L118").Value = 79999
$ws.Range("N118").Value = -83313
$ws.Range("H126").Value = 2762.5
$ws.Range("I126").Value = 2863.08
$ws.Range("J126").Value = 2259.6
$ws.Range("K126").Value = 8589.24
$ws.Range("L126").Value = 6778.799999999999
$ws.Range("M126").Value = -6119.24
$ws.Range("N126").Value = -11718.8
$ws.Range("H132").Value = 19208902
$ws.Range("I132").Value = 19208902
$ws.Range("K132").Value = 57626706
$ws.Range("M132").Value = -57624176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 87499
$ws.Range("J116").Value = 87499
$ws.Range("L116").Value = 87499
$ws.Range("N116").Value = -96677
$ws.Range("H126").Value = 2651.9092
$ws.Range("I126").Value = 2686.1667
$ws.Range("K126").Value = 8058.500100000001
$ws.Range("M126").Value = -5588.500100000001
$ws.Range("H132").Value = 17244878
$ws.Range("I132").Value = 25001766
$ws.Range("J132").Value = 7347.8887
$ws.Range("K132").Value = 75005298
$ws.Range("L132").Value = 22043.6661
$ws.Range("M132").Value = -75002768
$ws.Range("N132").Value = -27103.6661
